$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.160.59"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.33"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.29"
$ws.Range("E5").Value = "  +3.91%  "

$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.84"
$ws.Range("E8").Value = "  +6.57%  "

$ws.Range("E9").Value = "  +2.79%  "

$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.125.76"
$ws.Range("E12").Value = "  +2.00%  "

$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.850.76"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.141.57"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.67"

$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.38"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.23"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +2.50%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.76"

$ws.Range("E26").Value = "  +27.99%  "

$ws.Range("E27").Value = "  +3.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0558"
$ws.Range("E31").Value = "  +2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.99"
$ws.Range("E32").Value = "  +2.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +28.49%  "

$ws.Range("E34").Value = "  +2.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.829"
$ws.Range("E35").Value = "  +19.47%  "

$ws.Range("E36").Value = "  +10.22%  "

$ws.Range("E37").Value = "  +7.06%  "

$ws.Range("E38").Value = "  +7.73%  "

$ws.Range("E39").Value = "  +4.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "89.67"
$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.342.06"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.91"
$ws.Range("E42").Value = "  +3.93%  "

$ws.Range("E43").Value = "  +4.34%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0557"
$ws.Range("E45").Value = "  +6.79%  "

$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.58"
$ws.Range("E47").Value = "  +5.46%  "

$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.20"
$ws.Range("E48").Value = "  +44.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.039.78"
$ws.Range("E49").Value = "  +1.99%  "

$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("E51").Value = "  +0.19%  "
